$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C28").Value = 299
$ws.Range("D28").Value = 30
$ws.Range("E28").Value = 269
$ws.Range("F28").Value = 4.672897196261682
